$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.961.43"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "'2.243.24"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'271.05"
$ws.Range("E5").Value = "  +4.35%  "
$ws.Range("D6").Value = "'95.03"
$ws.Range("E6").Value = "  +15.72%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.635"
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("D10").Value = "'46.42"
$ws.Range("E10").Value = "  +7.49%  "
$ws.Range("D11").Value = "'0.0963"
$ws.Range("E11").Value = "  +5.05%  "
$ws.Range("D12").Value = "'8.27"
$ws.Range("E12").Value = "  +18.91%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "'15.30"
$ws.Range("E14").Value = "  +7.62%  "
$ws.Range("D15").Value = "'2.581.84"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("D17").Value = "'2.253.14"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "'43.935.33"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").Value = "'70.93"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").Value = "'235.45"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'9.18"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'11.44"
$ws.Range("E26").Value = "  +7.38%  "
$ws.Range("E27").Value = "  +12.60%  "
$ws.Range("D28").Value = "'3.59"
$ws.Range("E28").Value = "  +6.79%  "
$ws.Range("D29").Value = "'40.45"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'173.00"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'0.0915"
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("D33").Value = "'21.03"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("E34").Value = "  +3.66%  "
$ws.Range("D35").Value = "'0.125"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("E39").Value = "  +26.54%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'12.88"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.227"
$ws.Range("E41").Value = "  +13.81%  "
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("D43").Value = "'63.64"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "'0.0999"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "'101.17"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "'2.465.23"
$ws.Range("E51").Value = "  +2.27%  "
